$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("mysql iltrowania wynikow", $true, $false, $false, $false, $false, $true, 1, $false, "mysql filtrowania wynikow", 2)
